$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated capital structure database values for rows 2 and 3.
# Both rows receive identical updated figures across columns G:AQ.
$values = [ordered]@{
    "G"  = -0.7291666666666667
    "H"  = -0.8154761904761906
    "I"  = -0.755952380952381
    "J"  = -0.755952380952381
    "K"  = -1.09
    "L"  = -0.6488095238095238
    "M"  = 2.46
    "N"  = 0.06212121212121212
    "O"  = -2.256880733944954
    "S"  = 2.46
    "T"  = 1
    "U"  = 1.09
    "V"  = 0.02752525252525253
    "W"  = -0.6193181818181819
    "X"  = 0.06860776252520533
    "Y"  = -0.6879259443433872
    "Z"  = 1.228070175438597
    "AA" = -0.9283625730994153
    "AB" = 0.06806180419959879
    "AC" = -0.996424377299014
    "AD" = 0.451
    "AF" = 0.451
    "AG" = -0.639
    "AH" = 0.01126064268058226
    "AI" = 0.1560013836042892
    "AJ" = -0.0164010164010164
    "AK" = -0.3548028872848418
    "AL" = 0.032
    "AM" = 0.032
    "AN" = -0.4063063063063063
    "AO" = -39.6875
    "AP" = 0.5756756756756757
    "AQ" = -39.6875
}

foreach ($col in $values.Keys) {
    $val = $values[$col]
    $ws.Range("${col}2").Value = $val
    $ws.Range("${col}3").Value = $val
}
